# WAT new API's automation
# Adds 5 new test rows (29-33) to the WoS_AuthorTransformation sheet for the
# new "search author cluster by ORCID/RID" API endpoints.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlTop = -4160
$xlContinuous = 1
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10

# Column order used throughout the sheet:
# A=TESTNAME  B=DESCRIPTION  C=HOST  D=APIPATH  E=METHOD  F=HEADERS
# G=QUERYSTRING  H=BODY  I=DEPENDENCYTESTS  J=VALIDATIONS  K=STORE  L=STATUS

$rows = @(
    @{
        Row = 29
        Height = 120
        A = "WAT-408"
        B = "Verify that user is able to search for author cluster using ORCID"
        C = "1PRECOMMEND"
        D = "recommend/search/author/clusters/0000-0002-6423-7213"
        E = "GET"
        J = "status=200||hits.authorClusterId=24303705"
        K = "hits[0].authorClusterId||hits[0].primaryName||hits[0].alternativeNames||hits[0].affiliation||hits[0].location||hits[0].totalNumberOfPublications||hits[0].publicationYearRangeMin||hits[0].publicationYearRangeMax||hits[0].topPublications"
        KTop = $true
    },
    @{
        Row = 30
        Height = 120
        A = "WAT-504"
        B = "Verify that user is able to search for author cluster using RID"
        C = "1PRECOMMEND"
        D = "recommend/search/author/clusters/A-9832-2009"
        E = "GET"
        J = "status=200||hits.authorClusterId=80453160"
        K = "hits[0].authorClusterId||hits[0].primaryName||hits[0].alternativeNames||hits[0].affiliation||hits[0].location||hits[0].totalNumberOfPublications||hits[0].publicationYearRangeMin||hits[0].publicationYearRangeMax||hits[0].topPublications"
        KTop = $true
    },
    @{
        Row = 31
        Height = $null
        A = "WAT-409"
        B = "Verify that user is able to search author cluster when ORCID/RID is missing in the request"
        C = "1PRECOMMEND"
        D = "recommend/search/author/clusters/"
        E = "GET"
        J = "status=200||info.totalHits=0"
        K = $null
        KTop = $false
    },
    @{
        Row = 32
        Height = $null
        A = "WAT-505"
        B = "Verify that user is able to search author cluster using invalid ORCID"
        C = "1PRECOMMEND"
        D = "recommend/search/author/clusters/0000-abcd-6423-12ec"
        E = "GET"
        J = "status=200||info.totalHits=0"
        K = $null
        KTop = $false
    },
    @{
        Row = 33
        Height = $null
        A = "WAT-506"
        B = "Verify that user is able to search author cluster using invalid RID"
        C = "1PRECOMMEND"
        D = "recommend/search/author/clusters/A-1456-abcs"
        E = "GET"
        J = "status=200||info.totalHits=0"
        K = $null
        KTop = $false
    }
)

$xlNone = -4142

# The old last row (28) no longer sits at the bottom of the table, so the
# bottom border of its HEADERS..STORE cells (which used to close off the
# table) is dropped - matching the workbook's real "insert rows below"
# behaviour.
$ws.Range("F28:K28").Borders.Item($xlEdgeBottom).LineStyle = $xlNone

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Range("A$rowNum").Value = $r.A
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
    $ws.Range("E$rowNum").Value = $r.E
    $ws.Range("J$rowNum").Value = $r.J
    if ($r.K -ne $null) {
        $ws.Range("K$rowNum").Value = $r.K
    }

    # Full box border (thin) around every individual cell A:L, matching the
    # rest of the table (each cell - not just the outer edge of the range -
    # carries its own left/top/bottom/right border).
    for ($c = 1; $c -le 12; $c++) {
        $cell = $ws.Cells.Item($rowNum, $c)
        $cell.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
        $cell.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
        $cell.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
        $cell.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
    }

    # Vertical-top alignment for the narrower text columns.
    $ws.Range("A$rowNum").VerticalAlignment = $xlTop
    $ws.Range("C$rowNum").VerticalAlignment = $xlTop
    $ws.Range("D$rowNum").VerticalAlignment = $xlTop
    $ws.Range("E$rowNum").VerticalAlignment = $xlTop

    if ($r.KTop) {
        $ws.Range("K$rowNum").VerticalAlignment = $xlTop
        $ws.Range("K$rowNum").WrapText = $true
    }

    if ($r.Height -ne $null) {
        $ws.Rows.Item($rowNum).RowHeight = $r.Height
    }
}

# Match the workbook's final selection / scroll state from the edit.
$ws.Range("D33").Select()

Write-Output "Added WAT-408, WAT-504, WAT-409, WAT-505, WAT-506 rows (29-33)"
